# Updates market/profit data (columns H-N) on the per-job Leve tables,
# refreshed via the scheduled Diabolos Profits market-data runner.
$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 47035.875
$ws.Range("I33").Value = 62680.418
$ws.Range("J33").Value = 102.25
$ws.Range("K33").Value = 62680.418
$ws.Range("L33").Value = 102.25
$ws.Range("M33").Value = -62451.418
$ws.Range("N33").Value = -560.25
# Row 53
$ws.Range("H53").Value = 2884.818
$ws.Range("J53").Value = 4839.077
$ws.Range("L53").Value = 4839.077
$ws.Range("N53").Value = -6113.077
# Row 80
$ws.Range("H80").Value = 843260.4
$ws.Range("I80").Value = 1422094.2
$ws.Range("J80").Value = 1320.2727
$ws.Range("K80").Value = 4266282.6
$ws.Range("L80").Value = 3960.8181
$ws.Range("M80").Value = -4265284.6
$ws.Range("N80").Value = -5956.8181
# Row 83
$ws.Range("H83").Value = 843260.4
$ws.Range("I83").Value = 1422094.2
$ws.Range("J83").Value = 1320.2727
$ws.Range("K83").Value = 12798847.8
$ws.Range("L83").Value = 11882.4543
$ws.Range("M83").Value = -12793855.8
$ws.Range("N83").Value = -21866.4543
# Row 88
$ws.Range("H88").Value = 1898.8667
$ws.Range("I88").Value = 2062.25
$ws.Range("K88").Value = 2062.25
$ws.Range("M88").Value = -1656.25
# Row 91
$ws.Range("H91").Value = 1898.8667
$ws.Range("I91").Value = 2062.25
$ws.Range("K91").Value = 2062.25
$ws.Range("M91").Value = -658.25
# Row 96
$ws.Range("H96").Value = 1052.2
$ws.Range("I96").Value = 690.25
$ws.Range("K96").Value = 2070.75
$ws.Range("M96").Value = -697.75
# Row 103
$ws.Range("H103").Value = 511
$ws.Range("I103").Value = 316.75
$ws.Range("J103").Value = 770
$ws.Range("K103").Value = 950.25
$ws.Range("L103").Value = 2310
$ws.Range("M103").Value = -364.25
$ws.Range("N103").Value = -3482
# Row 132
$ws.Range("H132").Value = 3441.1462
$ws.Range("I132").Value = 2989.4358
$ws.Range("K132").Value = 8968.307400000002
$ws.Range("M132").Value = -6438.307400000002
# Row 137
$ws.Range("H137").Value = 3397.2
$ws.Range("I137").Value = 2990
$ws.Range("K137").Value = 8970
$ws.Range("M137").Value = -6420

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 88
$ws.Range("H88").Value = 10418413
$ws.Range("I88").Value = 27778678
$ws.Range("K88").Value = 27778678
$ws.Range("M88").Value = -27778272
# Row 91
$ws.Range("H91").Value = 10418413
$ws.Range("I91").Value = 27778678
$ws.Range("K91").Value = 27778678
$ws.Range("M91").Value = -27777274
# Row 102
$ws.Range("H102").Value = 12752.333
$ws.Range("I102").Value = 3731.3
$ws.Range("K102").Value = 3731.3
$ws.Range("M102").Value = -2109.3

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 82
$ws.Range("H82").Value = 41352.715
$ws.Range("I82").Value = 11167.833
$ws.Range("K82").Value = 11167.833
$ws.Range("M82").Value = -10784.833
# Row 85
$ws.Range("H85").Value = 41352.715
$ws.Range("I85").Value = 11167.833
$ws.Range("K85").Value = 11167.833
$ws.Range("M85").Value = -9841.833000000001
# Row 86
$ws.Range("H86").Value = 3501.375
$ws.Range("I86").Value = 2668.5
$ws.Range("J86").Value = 6000
$ws.Range("K86").Value = 2668.5
$ws.Range("L86").Value = 6000
$ws.Range("M86").Value = -1545.5
$ws.Range("N86").Value = -8246
# Row 89
$ws.Range("H89").Value = 3501.375
$ws.Range("I89").Value = 2668.5
$ws.Range("J89").Value = 6000
$ws.Range("K89").Value = 13342.5
$ws.Range("L89").Value = 30000
$ws.Range("M89").Value = -7726.5
$ws.Range("N89").Value = -41232
# Row 94
$ws.Range("H94").Value = 2818.2964
$ws.Range("I94").Value = 4037.2856
$ws.Range("K94").Value = 4037.2856
$ws.Range("M94").Value = -3586.2856
# Row 99
$ws.Range("H99").Value = 908.2
$ws.Range("I99").Value = 782.75
$ws.Range("K99").Value = 782.75
$ws.Range("M99").Value = 715.25

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2995.8984
$ws.Range("I31").Value = 1610.0588
$ws.Range("J31").Value = 3448.9614
$ws.Range("K31").Value = 1610.0588
$ws.Range("L31").Value = 3448.9614
$ws.Range("M31").Value = -1315.0588
$ws.Range("N31").Value = -4038.9614
# Row 34
$ws.Range("H34").Value = 2995.8984
$ws.Range("I34").Value = 1610.0588
$ws.Range("J34").Value = 3448.9614
$ws.Range("K34").Value = 1610.0588
$ws.Range("L34").Value = 3448.9614
$ws.Range("M34").Value = -1408.0588
$ws.Range("N34").Value = -3852.9614
# Row 86
$ws.Range("H86").Value = 4799.476
$ws.Range("I86").Value = 5015.933
$ws.Range("J86").Value = 4258.3335
$ws.Range("K86").Value = 5015.933
$ws.Range("L86").Value = 4258.3335
$ws.Range("M86").Value = -3892.933
$ws.Range("N86").Value = -6504.3335
# Row 89
$ws.Range("H89").Value = 4799.476
$ws.Range("I89").Value = 5015.933
$ws.Range("J89").Value = 4258.3335
$ws.Range("K89").Value = 25079.665
$ws.Range("L89").Value = 21291.6675
$ws.Range("M89").Value = -19463.665
$ws.Range("N89").Value = -32523.6675
# Row 105
$ws.Range("H105").Value = 2516.6667
$ws.Range("I105").Value = 1250
$ws.Range("J105").Value = 3150
$ws.Range("K105").Value = 1250
$ws.Range("L105").Value = 3150
$ws.Range("M105").Value = 497
$ws.Range("N105").Value = -6644

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 121
$ws.Range("H121").Value = 114225
$ws.Range("I121").Value = 2096.6667
$ws.Range("J121").Value = 156273.12
$ws.Range("K121").Value = 6290.000100000001
$ws.Range("L121").Value = 468819.36
$ws.Range("M121").Value = -4980.000100000001
$ws.Range("N121").Value = -471439.36

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 141
$ws.Range("H141").Value = 92249.75
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 92249.75
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 92249.75
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -102609.75

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 12502853
$ws.Range("I7").Value = 16669036
$ws.Range("K7").Value = 16669036
$ws.Range("M7").Value = -16668924
# Row 93
$ws.Range("H93").Value = 5740.4
$ws.Range("I93").Value = 5740.4
$ws.Range("K93").Value = 5740.4
$ws.Range("M93").Value = -4492.4
# Row 126
$ws.Range("H126").Value = 12502853
$ws.Range("I126").Value = 16669036
$ws.Range("K126").Value = 50007108
$ws.Range("M126").Value = -50004638

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 3315.8333
$ws.Range("I96").Value = 3315.8333
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 3315.8333
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -1942.8333
$ws.Range("N96").ClearContents()
# Row 100
$ws.Range("H100").Value = 5284.8
$ws.Range("I100").Value = 5838.6665
$ws.Range("K100").Value = 11677.333
$ws.Range("M100").Value = -11136.333
# Row 122
$ws.Range("H122").Value = 1852.6471
$ws.Range("I122").Value = 1806.7142
$ws.Range("K122").Value = 5420.142599999999
$ws.Range("M122").Value = -2970.142599999999
